# base_de_datos.xlsx update: append a new data row (row 9) that duplicates
# row 8's person ("Jose Luis" Roque, age 25, maria@yahoo.com) with the next
# sequential ID (8), then leave the view scrolled/selected the way the
# author left it (selection on D20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the new row of data -------------------------------------------
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "Jose Luis"
$ws.Cells.Item(9, 3).Value = "Roque"
$ws.Cells.Item(9, 4).Value = 25
$ws.Cells.Item(9, 5).Value = "maria@yahoo.com"

# --- Restore the window/view state recorded in the saved workbook ---------
# Author had scrolled the grid down (topLeftCell = A9) before saving.
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 1
$win.TabRatio = 0.58

# Final selection left on D20, as captured in the workbook's sheetView.
$ws.Range("D20").Select()
